$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-formatted columns (Coin name/link/price/volume change) to stay as text
# so purely-numeric-looking price strings (e.g. trailing-zero decimals) are not
# silently coerced into Number cells by Excel's automatic type inference.
$ws.Range('B2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '60.656.79'
$ws.Range('E2').Value = '  +2.18%  '
$ws.Range('D3').Value = '2.608.49'
$ws.Range('E3').Value = '  +1.14%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '578.62'
$ws.Range('E5').Value = '  +4.09%  '
$ws.Range('D6').Value = '143.16'
$ws.Range('E6').Value = '  +0.63%  '
$ws.Range('E7').Value = '  -0.27%  '
$ws.Range('E8').Value = '  +0.45%  '
$ws.Range('D9').Value = '2.634.12'
$ws.Range('E9').Value = '  +1.83%  '
$ws.Range('D10').Value = '6.53'
$ws.Range('E10').Value = '  -2.08%  '
$ws.Range('E11').Value = '  +1.86%  '
$ws.Range('D12').Value = '0.154'
$ws.Range('E12').Value = '  -6.59%  '
$ws.Range('E13').Value = '  +3.70%  '
$ws.Range('D14').Value = '3.071.84'
$ws.Range('E14').Value = '  +1.22%  '
$ws.Range('D15').Value = '60.656.36'
$ws.Range('E15').Value = '  +2.18%  '
$ws.Range('D16').Value = '23.28'
$ws.Range('E16').Value = '  +1.11%  '
$ws.Range('E17').Value = '  +3.97%  '
$ws.Range('D18').Value = '2.623.91'
$ws.Range('E18').Value = '  +1.45%  '
$ws.Range('D19').Value = '11.32'
$ws.Range('E19').Value = '  +9.55%  '
$ws.Range('D20').Value = '4.65'
$ws.Range('E20').Value = '  +2.08%  '
$ws.Range('D21').Value = '349.00'
$ws.Range('E21').Value = '  +3.46%  '
$ws.Range('D22').Value = '6.93'
$ws.Range('E22').Value = '  +7.82%  '
$ws.Range('E23').Value = '  -0.18%  '
$ws.Range('E24').Value = '  +11.22%  '
$ws.Range('D25').Value = '63.40'
$ws.Range('E25').Value = '  +0.39%  '
$ws.Range('E26').Value = '  -0.21%  '
$ws.Range('E27').Value = '  +0.42%  '
$ws.Range('E28').Value = '  +4.32%  '
$ws.Range('D29').Value = '0.0₃0795'
$ws.Range('E29').Value = '  +2.37%  '
$ws.Range('E30').Value = '  +12.28%  '
$ws.Range('E31').Value = '  +3.34%  '
$ws.Range('E32').Value = '  -0.09%  '
$ws.Range('D33').Value = '162.18'
$ws.Range('E33').Value = '  +2.77%  '
$ws.Range('D34').Value = '19.57'
$ws.Range('E34').Value = '  +2.56%  '
$ws.Range('E35').Value = '  +4.58%  '
$ws.Range('D36').Value = '0.973'
$ws.Range('E36').Value = '  +7.14%  '
$ws.Range('D37').Value = '1.24'
$ws.Range('E37').Value = '  +7.05%  '
$ws.Range('E38').Value = '  +7.92%  '
$ws.Range('D40').Value = '3.87'
$ws.Range('E40').Value = '  +5.37%  '
$ws.Range('E41').Value = '  -1.77%  '
$ws.Range('D42').Value = '300.57'
$ws.Range('E42').Value = '  +3.05%  '
$ws.Range('D43').Value = '135.10'
$ws.Range('E43').Value = '  -1.05%  '
$ws.Range('D44').Value = '0.996'
$ws.Range('E44').Value = '  -0.36%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '19.94'
$ws.Range('E45').Value = '  +4.92%  '
$ws.Range('B46').Value = 'Stellar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D46').Value = '0.0986'
$ws.Range('E46').Value = '  +1.10%  '
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').Value = '0.606'
$ws.Range('E47').Value = '  +2.28%  '
$ws.Range('B48').Value = 'Hedera'
$ws.Range('C48').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D48').Value = '0.0548'
$ws.Range('E48').Value = '  +3.00%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').Value = '4.98'
$ws.Range('E49').Value = '  +9.47%  '
$ws.Range('D50').Value = '0.0242'
$ws.Range('E50').Value = '  +3.30%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').Value = '19.85'
$ws.Range('E51').Value = '  +5.91%  '
